# Apply the "Generate Report for Archive" edit:
#  1. Update status text from "Ready for handoff" to "In Translation" wherever it appears
#     (Overview sheet columns zh-cn/de-de, and the per-language sheets' Status column).
#  2. Narrow the "Latest Handoff Datetime"-style date columns (Overview E/F and the
#     per-language sheets' Status column width companion) from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldStatus, $newStatus) | Out-Null
}

# Column width updates (closest value reachable through ColumnWidth's pixel snapping).
$targetColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
